# Daily "cryptos" price/volume refresh (GitHub Actions scrape update).
# All Price (D) / Volume(1h) (E) cells hold text, not numbers (e.g. "65.599.63",
# "2.952.90" use dots as thousands separators so they aren't valid numbers, and
# even plain-looking values like "0.999" must stay text to match the sheet's
# existing inline-string formatting). For any new D-column value that parses as
# a real number we force the cell to Text format first so Excel doesn't
# silently convert it to a numeric value.
# Rows 27/28, 36/37 and 41/42 also swap which coin occupies which row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.599.63"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.952.90"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.49"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.64"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "2.951.27"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.67"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "65.493.59"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "3.440.15"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.02"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "2.953.73"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.84"
$ws.Range("E20").Value = "  +13.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.48"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.26"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.29"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.49"
$ws.Range("E29").Value = "  +6.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.06"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.20"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.74"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.970"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "45.15"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.15"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -7.47%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.303"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("E43").Value = "  -5.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.54"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "386.83"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0352"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "2.681.96"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.02"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.83"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("E51").Value = "  +0.62%  "
